# Applies the row-shuffle + minor odds updates described in the commit diff.
# The underlying source data for several same-date fixtures was re-shuffled;
# only columns B (match id) and F..AC (teams/scores/odds) move between rows -
# column A (row index), C, D (Div / Div Original Name) and E (Date) stay fixed
# to their row position because every row in a block shares the same date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 576 <- original row 577
$ws.Range("B576").Value = 5207403
$ws.Range("F576").Value = '1860 Munich'
$ws.Range("G576").Value = 'VfB Oldenburg'
$ws.Range("H576").Value = 1
$ws.Range("I576").Value = 0
$ws.Range("J576").Value = 'H'
$ws.Range("K576").Value = 1.444
$ws.Range("L576").Value = 4.75
$ws.Range("M576").Value = 5.5
$ws.Range("N576").Value = 1.666
$ws.Range("O576").Value = 4.2
$ws.Range("P576").Value = 4
$ws.Range("Q576").Value = -0.75
$ws.Range("R576").Value = 1.85
$ws.Range("S576").Value = 1.95
$ws.Range("T576").Value = 3
$ws.Range("U576").Value = 2
$ws.Range("V576").Value = 1.8
$ws.Range("W576").Value = 0.6659999999999999
$ws.Range("X576").Value = -1
$ws.Range("Y576").Value = -1
$ws.Range("Z576").Value = 0.425
$ws.Range("AA576").Value = -0.5
$ws.Range("AB576").Value = -1
$ws.Range("AC576").Value = 0.8

# Row 577 <- original row 578
$ws.Range("B577").Value = 5207404
$ws.Range("F577").Value = 'Elversberg'
$ws.Range("G577").Value = 'Saarbrucken'
$ws.Range("H577").Value = 0
$ws.Range("I577").Value = 2
$ws.Range("J577").Value = 'A'
$ws.Range("K577").Value = 2.55
$ws.Range("L577").Value = 3.3
$ws.Range("M577").Value = 2.55
$ws.Range("N577").Value = 2.1
$ws.Range("O577").Value = 3.4
$ws.Range("P577").Value = 3.1
$ws.Range("Q577").Value = -0.25
$ws.Range("R577").Value = 1.8
$ws.Range("S577").Value = 2
$ws.Range("T577").Value = 2.75
$ws.Range("U577").Value = 1.8
$ws.Range("V577").Value = 2
$ws.Range("W577").Value = -1
$ws.Range("X577").Value = -1
$ws.Range("Y577").Value = 2.1
$ws.Range("Z577").Value = -1
$ws.Range("AA577").Value = 1
$ws.Range("AB577").Value = -1
$ws.Range("AC577").Value = 1

# Row 578 <- original row 579
$ws.Range("B578").Value = 5207932
$ws.Range("F578").Value = 'Hallescher FC'
$ws.Range("G578").Value = 'Dynamo Dresden'
$ws.Range("H578").Value = 0
$ws.Range("I578").Value = 2
$ws.Range("J578").Value = 'A'
$ws.Range("K578").Value = 3.5
$ws.Range("L578").Value = 3.5
$ws.Range("M578").Value = 1.909
$ws.Range("N578").Value = 3.3
$ws.Range("O578").Value = 3.3
$ws.Range("P578").Value = 2
$ws.Range("Q578").Value = 0.25
$ws.Range("R578").Value = 2.025
$ws.Range("S578").Value = 1.775
$ws.Range("T578").Value = 2.5
$ws.Range("U578").Value = 1.825
$ws.Range("V578").Value = 1.975
$ws.Range("W578").Value = -1
$ws.Range("X578").Value = -1
$ws.Range("Y578").Value = 1
$ws.Range("Z578").Value = -1
$ws.Range("AA578").Value = 0.7749999999999999
$ws.Range("AB578").Value = -1
$ws.Range("AC578").Value = 0.9750000000000001

# Row 579 <- original row 576
$ws.Range("B579").Value = 5212739
$ws.Range("F579").Value = 'FC Viktoria Kln'
$ws.Range("G579").Value = 'Wehen SV'
$ws.Range("H579").Value = 1
$ws.Range("I579").Value = 0
$ws.Range("J579").Value = 'H'
$ws.Range("K579").Value = 2.4
$ws.Range("L579").Value = 3.3
$ws.Range("M579").Value = 2.7
$ws.Range("N579").Value = 2.7
$ws.Range("O579").Value = 3.3
$ws.Range("P579").Value = 2.4
$ws.Range("Q579").Value = 0
$ws.Range("R579").Value = 2.025
$ws.Range("S579").Value = 1.775
$ws.Range("T579").Value = 2.5
$ws.Range("U579").Value = 1.825
$ws.Range("V579").Value = 1.975
$ws.Range("W579").Value = 1.7
$ws.Range("X579").Value = -1
$ws.Range("Y579").Value = -1
$ws.Range("Z579").Value = 1.025
$ws.Range("AA579").Value = -1
$ws.Range("AB579").Value = -1
$ws.Range("AC579").Value = 0.9750000000000001

# Row 639 <- original row 640
$ws.Range("B639").Value = 5207253
$ws.Range("F639").Value = 'SpVgg Bayreuth'
$ws.Range("G639").Value = 'Saarbrucken'
$ws.Range("H639").Value = 0
$ws.Range("I639").Value = 6
$ws.Range("J639").Value = 'A'
$ws.Range("K639").Value = 3.6
$ws.Range("L639").Value = 3.5
$ws.Range("M639").Value = 1.85
$ws.Range("N639").Value = 3.8
$ws.Range("O639").Value = 3.75
$ws.Range("P639").Value = 1.8
$ws.Range("Q639").Value = 0.5
$ws.Range("R639").Value = 2
$ws.Range("S639").Value = 1.8
$ws.Range("T639").Value = 3
$ws.Range("U639").Value = 2
$ws.Range("V639").Value = 1.8
$ws.Range("W639").Value = -1
$ws.Range("X639").Value = -1
$ws.Range("Y639").Value = 0.8
$ws.Range("Z639").Value = -1
$ws.Range("AA639").Value = 0.8
$ws.Range("AB639").Value = 1
$ws.Range("AC639").Value = -1

# Row 640 <- original row 639
$ws.Range("B640").Value = 5207420
$ws.Range("F640").Value = 'Borussia Dortmund II'
$ws.Range("G640").Value = 'VfB Oldenburg'
$ws.Range("H640").Value = 1
$ws.Range("I640").Value = 2
$ws.Range("J640").Value = 'A'
$ws.Range("K640").Value = 1.95
$ws.Range("L640").Value = 3.5
$ws.Range("M640").Value = 3.5
$ws.Range("N640").Value = 2.25
$ws.Range("O640").Value = 3.3
$ws.Range("P640").Value = 2.875
$ws.Range("Q640").Value = -0.25
$ws.Range("R640").Value = 2
$ws.Range("S640").Value = 1.8
$ws.Range("T640").Value = 2.75
$ws.Range("U640").Value = 2
$ws.Range("V640").Value = 1.8
$ws.Range("W640").Value = -1
$ws.Range("X640").Value = -1
$ws.Range("Y640").Value = 1.875
$ws.Range("Z640").Value = -1
$ws.Range("AA640").Value = 0.8
$ws.Range("AB640").Value = 0.5
$ws.Range("AC640").Value = -0.5

# Row 747 <- original row 748
$ws.Range("B747").Value = 5207242
$ws.Range("F747").Value = 'Erzgebirge Aue'
$ws.Range("G747").Value = 'SpVgg Bayreuth'
$ws.Range("H747").Value = 4
$ws.Range("I747").Value = 0
$ws.Range("J747").Value = 'H'
$ws.Range("K747").Value = 1.727
$ws.Range("L747").Value = 3.6
$ws.Range("M747").Value = 4.2
$ws.Range("N747").Value = 1.7
$ws.Range("O747").Value = 3.75
$ws.Range("P747").Value = 4.333
$ws.Range("Q747").Value = -0.75
$ws.Range("R747").Value = 1.95
$ws.Range("S747").Value = 1.85
$ws.Range("T747").Value = 2.5
$ws.Range("U747").Value = 1.8
$ws.Range("V747").Value = 2
$ws.Range("W747").Value = 0.7
$ws.Range("X747").Value = -1
$ws.Range("Y747").Value = -1
$ws.Range("Z747").Value = 0.95
$ws.Range("AA747").Value = -1
$ws.Range("AB747").Value = 0.8
$ws.Range("AC747").Value = -1

# Row 748 <- original row 749
$ws.Range("B748").Value = 5207449
$ws.Range("F748").Value = 'Elversberg'
$ws.Range("G748").Value = 'FC Ingolstadt'
$ws.Range("H748").Value = 4
$ws.Range("I748").Value = 3
$ws.Range("J748").Value = 'H'
$ws.Range("K748").Value = 2.1
$ws.Range("L748").Value = 3.4
$ws.Range("M748").Value = 3.1
$ws.Range("N748").Value = 1.8
$ws.Range("O748").Value = 3.6
$ws.Range("P748").Value = 3.8
$ws.Range("Q748").Value = -0.5
$ws.Range("R748").Value = 1.825
$ws.Range("S748").Value = 1.975
$ws.Range("T748").Value = 2.75
$ws.Range("U748").Value = 1.9
$ws.Range("V748").Value = 1.9
$ws.Range("W748").Value = 0.8
$ws.Range("X748").Value = -1
$ws.Range("Y748").Value = -1
$ws.Range("Z748").Value = 0.825
$ws.Range("AA748").Value = -1
$ws.Range("AB748").Value = 0.8999999999999999
$ws.Range("AC748").Value = -1

# Row 749 <- original row 747
$ws.Range("B749").Value = 5207450
$ws.Range("F749").Value = 'Verl'
$ws.Range("G749").Value = 'RotWeiss Essen'
$ws.Range("H749").Value = 1
$ws.Range("I749").Value = 1
$ws.Range("J749").Value = 'D'
$ws.Range("K749").Value = 2.25
$ws.Range("L749").Value = 3.3
$ws.Range("M749").Value = 2.9
$ws.Range("N749").Value = 2.2
$ws.Range("O749").Value = 3.3
$ws.Range("P749").Value = 3
$ws.Range("Q749").Value = -0.25
$ws.Range("R749").Value = 1.925
$ws.Range("S749").Value = 1.875
$ws.Range("T749").Value = 2.5
$ws.Range("U749").Value = 1.825
$ws.Range("V749").Value = 1.975
$ws.Range("W749").Value = -1
$ws.Range("X749").Value = 2.3
$ws.Range("Y749").Value = -1
$ws.Range("Z749").Value = -0.5
$ws.Range("AA749").Value = 0.4375
$ws.Range("AB749").Value = -1
$ws.Range("AC749").Value = 0.9750000000000001

# Row 776 <- original row 778
$ws.Range("B776").Value = 5207456
$ws.Range("F776").Value = 'Elversberg'
$ws.Range("G776").Value = 'VfB Oldenburg'
$ws.Range("H776").Value = 3
$ws.Range("I776").Value = 0
$ws.Range("J776").Value = 'H'
$ws.Range("K776").Value = 1.363
$ws.Range("L776").Value = 4.75
$ws.Range("M776").Value = 6.5
$ws.Range("N776").Value = 1.363
$ws.Range("O776").Value = 5
$ws.Range("P776").Value = 6
$ws.Range("Q776").Value = -1.25
$ws.Range("R776").Value = 1.8
$ws.Range("S776").Value = 2
$ws.Range("T776").Value = 3.25
$ws.Range("U776").Value = 1.975
$ws.Range("V776").Value = 1.825
$ws.Range("W776").Value = 0.363
$ws.Range("X776").Value = -1
$ws.Range("Y776").Value = -1
$ws.Range("Z776").Value = 0.8
$ws.Range("AA776").Value = -1
$ws.Range("AB776").Value = -0.5
$ws.Range("AC776").Value = 0.4125

# Row 777 <- original row 779
$ws.Range("B777").Value = 5207948
$ws.Range("F777").Value = 'Verl'
$ws.Range("G777").Value = 'Dynamo Dresden'
$ws.Range("H777").Value = 2
$ws.Range("I777").Value = 3
$ws.Range("J777").Value = 'A'
$ws.Range("K777").Value = 2.45
$ws.Range("L777").Value = 3.25
$ws.Range("M777").Value = 2.625
$ws.Range("N777").Value = 3
$ws.Range("O777").Value = 3.2
$ws.Range("P777").Value = 2.2
$ws.Range("Q777").Value = 0.25
$ws.Range("R777").Value = 1.825
$ws.Range("S777").Value = 1.975
$ws.Range("T777").Value = 2.5
$ws.Range("U777").Value = 1.85
$ws.Range("V777").Value = 1.95
$ws.Range("W777").Value = -1
$ws.Range("X777").Value = -1
$ws.Range("Y777").Value = 1.2
$ws.Range("Z777").Value = -1
$ws.Range("AA777").Value = 0.9750000000000001
$ws.Range("AB777").Value = 0.8500000000000001
$ws.Range("AC777").Value = -1

# Row 778 <- original row 777
$ws.Range("B778").Value = 5212493
$ws.Range("F778").Value = 'SV Meppen'
$ws.Range("G778").Value = '1860 Munich'
$ws.Range("H778").Value = 2
$ws.Range("I778").Value = 1
$ws.Range("J778").Value = 'H'
$ws.Range("K778").Value = 3.1
$ws.Range("L778").Value = 3.4
$ws.Range("M778").Value = 2.1
$ws.Range("N778").Value = 3.2
$ws.Range("O778").Value = 3.4
$ws.Range("P778").Value = 2.05
$ws.Range("Q778").Value = 0.25
$ws.Range("R778").Value = 1.95
$ws.Range("S778").Value = 1.85
$ws.Range("T778").Value = 2.75
$ws.Range("U778").Value = 1.85
$ws.Range("V778").Value = 1.95
$ws.Range("W778").Value = 2.2
$ws.Range("X778").Value = -1
$ws.Range("Y778").Value = -1
$ws.Range("Z778").Value = 0.95
$ws.Range("AA778").Value = -1
$ws.Range("AB778").Value = 0.425
$ws.Range("AC778").Value = -0.5

# Row 779 <- original row 776
$ws.Range("B779").Value = 5212703
$ws.Range("F779").Value = 'Borussia Dortmund II'
$ws.Range("G779").Value = 'Saarbrucken'
$ws.Range("H779").Value = 1
$ws.Range("I779").Value = 2
$ws.Range("J779").Value = 'A'
$ws.Range("K779").Value = 3.2
$ws.Range("L779").Value = 3.3
$ws.Range("M779").Value = 2.05
$ws.Range("N779").Value = 3
$ws.Range("O779").Value = 3.2
$ws.Range("P779").Value = 2.2
$ws.Range("Q779").Value = 0.25
$ws.Range("R779").Value = 1.875
$ws.Range("S779").Value = 1.925
$ws.Range("T779").Value = 2.5
$ws.Range("U779").Value = 1.9
$ws.Range("V779").Value = 1.9
$ws.Range("W779").Value = -1
$ws.Range("X779").Value = -1
$ws.Range("Y779").Value = 1.2
$ws.Range("Z779").Value = -1
$ws.Range("AA779").Value = 0.925
$ws.Range("AB779").Value = 0.8999999999999999
$ws.Range("AC779").Value = -1

# Row 886 <- original row 889
$ws.Range("B886").Value = 5447797
$ws.Range("F886").Value = 'MSV Duisburg'
$ws.Range("G886").Value = 'Wehen SV'
$ws.Range("H886").Value = 1
$ws.Range("I886").Value = 1
$ws.Range("J886").Value = 'D'
$ws.Range("K886").Value = 3.5
$ws.Range("L886").Value = 3.5
$ws.Range("M886").Value = 1.85
$ws.Range("N886").Value = 4.333
$ws.Range("O886").Value = 3.75
$ws.Range("P886").Value = 1.7
$ws.Range("Q886").Value = 0.5
$ws.Range("R886").Value = 2.2
$ws.Range("S886").Value = 1.7
$ws.Range("T886").Value = 2.75
$ws.Range("U886").Value = 1.825
$ws.Range("V886").Value = 2.025
$ws.Range("W886").Value = -1
$ws.Range("X886").Value = 2.75
$ws.Range("Y886").Value = -1
$ws.Range("Z886").Value = 1.2
$ws.Range("AA886").Value = -1
$ws.Range("AB886").Value = -1
$ws.Range("AC886").Value = 1.025

# Row 889 <- original row 886
$ws.Range("B889").Value = 5447800
$ws.Range("F889").Value = '1860 Munich'
$ws.Range("G889").Value = 'SpVgg Bayreuth'
$ws.Range("H889").Value = 2
$ws.Range("I889").Value = 0
$ws.Range("J889").Value = 'H'
$ws.Range("K889").Value = 1.5
$ws.Range("L889").Value = 4.2
$ws.Range("M889").Value = 5
$ws.Range("N889").Value = 1.533
$ws.Range("O889").Value = 4.333
$ws.Range("P889").Value = 4.75
$ws.Range("Q889").Value = -1
$ws.Range("R889").Value = 1.875
$ws.Range("S889").Value = 1.925
$ws.Range("T889").Value = 3.25
$ws.Range("U889").Value = 2
$ws.Range("V889").Value = 1.8
$ws.Range("W889").Value = 0.5329999999999999
$ws.Range("X889").Value = -1
$ws.Range("Y889").Value = -1
$ws.Range("Z889").Value = 0.875
$ws.Range("AA889").Value = -1
$ws.Range("AB889").Value = -1
$ws.Range("AC889").Value = 0.8

# Row 896 <- original row 898
$ws.Range("B896").Value = 5465120
$ws.Range("F896").Value = 'Borussia Dortmund II'
$ws.Range("G896").Value = 'Elversberg'
$ws.Range("H896").Value = 2
$ws.Range("I896").Value = 0
$ws.Range("J896").Value = 'H'
$ws.Range("K896").Value = 3.8
$ws.Range("L896").Value = 3.75
$ws.Range("M896").Value = 1.833
$ws.Range("N896").Value = 4
$ws.Range("O896").Value = 3.8
$ws.Range("P896").Value = 1.8
$ws.Range("Q896").Value = 0.5
$ws.Range("R896").Value = 2
$ws.Range("S896").Value = 1.85
$ws.Range("T896").Value = 3
$ws.Range("U896").Value = 2.025
$ws.Range("V896").Value = 1.825
$ws.Range("W896").Value = 3
$ws.Range("X896").Value = -1
$ws.Range("Y896").Value = -1
$ws.Range("Z896").Value = 1
$ws.Range("AA896").Value = -1
$ws.Range("AB896").Value = -1
$ws.Range("AC896").Value = 0.825

# Row 897 <- original row 896
$ws.Range("B897").Value = 5465166
$ws.Range("F897").Value = 'Freiburg II'
$ws.Range("G897").Value = 'Dynamo Dresden'
$ws.Range("H897").Value = 1
$ws.Range("I897").Value = 1
$ws.Range("J897").Value = 'D'
$ws.Range("K897").Value = 2.3
$ws.Range("L897").Value = 3.25
$ws.Range("M897").Value = 3
$ws.Range("N897").Value = 2.375
$ws.Range("O897").Value = 3.25
$ws.Range("P897").Value = 2.9
$ws.Range("Q897").Value = -0.25
$ws.Range("R897").Value = 2.05
$ws.Range("S897").Value = 1.8
$ws.Range("T897").Value = 2.5
$ws.Range("U897").Value = 1.925
$ws.Range("V897").Value = 1.925
$ws.Range("W897").Value = -1
$ws.Range("X897").Value = 2.25
$ws.Range("Y897").Value = -1
$ws.Range("Z897").Value = -0.5
$ws.Range("AA897").Value = 0.4
$ws.Range("AB897").Value = -1
$ws.Range("AC897").Value = 0.925

# Row 898 <- original row 899
$ws.Range("B898").Value = 5470244
$ws.Range("F898").Value = 'Erzgebirge Aue'
$ws.Range("G898").Value = 'FC Viktoria Kln'
$ws.Range("H898").Value = 1
$ws.Range("I898").Value = 1
$ws.Range("J898").Value = 'D'
$ws.Range("K898").Value = 2.5
$ws.Range("L898").Value = 3.3
$ws.Range("M898").Value = 2.7
$ws.Range("N898").Value = 2.875
$ws.Range("O898").Value = 3.4
$ws.Range("P898").Value = 2.3
$ws.Range("Q898").Value = 0.25
$ws.Range("R898").Value = 1.775
$ws.Range("S898").Value = 2.025
$ws.Range("T898").Value = 3
$ws.Range("U898").Value = 2
$ws.Range("V898").Value = 1.8
$ws.Range("W898").Value = -1
$ws.Range("X898").Value = 2.4
$ws.Range("Y898").Value = -1
$ws.Range("Z898").Value = 0.3875
$ws.Range("AA898").Value = -0.5
$ws.Range("AB898").Value = -1
$ws.Range("AC898").Value = 0.8

# Row 899 <- original row 897
$ws.Range("B899").Value = 5465167
$ws.Range("F899").Value = 'SpVgg Bayreuth'
$ws.Range("G899").Value = 'MSV Duisburg'
$ws.Range("H899").Value = 0
$ws.Range("I899").Value = 4
$ws.Range("J899").Value = 'A'
$ws.Range("K899").Value = 2.625
$ws.Range("L899").Value = 3.3
$ws.Range("M899").Value = 2.45
$ws.Range("N899").Value = 2.875
$ws.Range("O899").Value = 3.3
$ws.Range("P899").Value = 2.3
$ws.Range("Q899").Value = 0.25
$ws.Range("R899").Value = 1.825
$ws.Range("S899").Value = 2.025
$ws.Range("T899").Value = 3
$ws.Range("U899").Value = 2
$ws.Range("V899").Value = 1.85
$ws.Range("W899").Value = -1
$ws.Range("X899").Value = -1
$ws.Range("Y899").Value = 1.3
$ws.Range("Z899").Value = -1
$ws.Range("AA899").Value = 1.025
$ws.Range("AB899").Value = 1
$ws.Range("AC899").Value = -1

# Row 1025 <- original row 1027
$ws.Range("B1025").Value = 6881623
$ws.Range("F1025").Value = 'Saarbrucken'
$ws.Range("G1025").Value = 'Vfb Lubeck'
$ws.Range("H1025").Value = 1
$ws.Range("I1025").Value = 1
$ws.Range("J1025").Value = 'D'
$ws.Range("K1025").Value = 1.65
$ws.Range("L1025").Value = 3.6
$ws.Range("M1025").Value = 4.333
$ws.Range("N1025").Value = 1.5
$ws.Range("O1025").Value = 4
$ws.Range("P1025").Value = 5.75
$ws.Range("Q1025").Value = -1
$ws.Range("R1025").Value = 1.825
$ws.Range("S1025").Value = 1.975
$ws.Range("T1025").Value = 3
$ws.Range("U1025").Value = 1.875
$ws.Range("V1025").Value = 1.925
$ws.Range("W1025").Value = -1
$ws.Range("X1025").Value = 3
$ws.Range("Y1025").Value = -1
$ws.Range("Z1025").Value = -1
$ws.Range("AA1025").Value = 0.9750000000000001
$ws.Range("AB1025").Value = -1
$ws.Range("AC1025").Value = 0.925

# Row 1026 <- original row 1025
$ws.Range("B1026").Value = 6881624
$ws.Range("F1026").Value = 'Verl'
$ws.Range("G1026").Value = 'Arminia Bielefeld'
$ws.Range("H1026").Value = 3
$ws.Range("I1026").Value = 1
$ws.Range("J1026").Value = 'H'
$ws.Range("K1026").Value = 2.15
$ws.Range("L1026").Value = 3.3
$ws.Range("M1026").Value = 3.1
$ws.Range("N1026").Value = 2.1
$ws.Range("O1026").Value = 3.3
$ws.Range("P1026").Value = 3.2
$ws.Range("Q1026").Value = -0.25
$ws.Range("R1026").Value = 1.875
$ws.Range("S1026").Value = 1.925
$ws.Range("T1026").Value = 3
$ws.Range("U1026").Value = 1.95
$ws.Range("V1026").Value = 1.85
$ws.Range("W1026").Value = 1.1
$ws.Range("X1026").Value = -1
$ws.Range("Y1026").Value = -1
$ws.Range("Z1026").Value = 0.875
$ws.Range("AA1026").Value = -1
$ws.Range("AB1026").Value = 0.95
$ws.Range("AC1026").Value = -1

# Row 1027 <- original row 1026
$ws.Range("B1027").Value = 6881321
$ws.Range("F1027").Value = 'Freiburg II'
$ws.Range("G1027").Value = 'Sandhausen'
$ws.Range("H1027").Value = 0
$ws.Range("I1027").Value = 2
$ws.Range("J1027").Value = 'A'
$ws.Range("K1027").Value = 2.875
$ws.Range("L1027").Value = 3.4
$ws.Range("M1027").Value = 2.2
$ws.Range("N1027").Value = 2.9
$ws.Range("O1027").Value = 3.6
$ws.Range("P1027").Value = 2.1
$ws.Range("Q1027").Value = 0.25
$ws.Range("R1027").Value = 1.9
$ws.Range("S1027").Value = 1.9
$ws.Range("T1027").Value = 2.75
$ws.Range("U1027").Value = 1.925
$ws.Range("V1027").Value = 1.875
$ws.Range("W1027").Value = -1
$ws.Range("X1027").Value = -1
$ws.Range("Y1027").Value = 1.1
$ws.Range("Z1027").Value = -1
$ws.Range("AA1027").Value = 0.8999999999999999
$ws.Range("AB1027").Value = -1
$ws.Range("AC1027").Value = 0.875

# Row 1028 <- original row 1029
$ws.Range("B1028").Value = 6881322
$ws.Range("F1028").Value = 'SSV Ulm 1846'
$ws.Range("G1028").Value = '1860 Munich'
$ws.Range("H1028").Value = 1
$ws.Range("I1028").Value = 0
$ws.Range("J1028").Value = 'H'
$ws.Range("K1028").Value = 2.25
$ws.Range("L1028").Value = 3.4
$ws.Range("M1028").Value = 2.875
$ws.Range("N1028").Value = 2.3
$ws.Range("O1028").Value = 3.3
$ws.Range("P1028").Value = 2.75
$ws.Range("Q1028").Value = 0
$ws.Range("R1028").Value = 1.75
$ws.Range("S1028").Value = 2.05
$ws.Range("T1028").Value = 2.5
$ws.Range("U1028").Value = 1.9
$ws.Range("V1028").Value = 1.9
$ws.Range("W1028").Value = 1.3
$ws.Range("X1028").Value = -1
$ws.Range("Y1028").Value = -1
$ws.Range("Z1028").Value = 0.75
$ws.Range("AA1028").Value = -1
$ws.Range("AB1028").Value = -1
$ws.Range("AC1028").Value = 0.8999999999999999

# Row 1029 <- original row 1028
$ws.Range("B1029").Value = 6880358
$ws.Range("F1029").Value = 'FC Viktoria Kln'
$ws.Range("G1029").Value = 'FC Ingolstadt'
$ws.Range("H1029").Value = 1
$ws.Range("I1029").Value = 0
$ws.Range("J1029").Value = 'H'
$ws.Range("K1029").Value = 2.4
$ws.Range("L1029").Value = 3.5
$ws.Range("M1029").Value = 2.5
$ws.Range("N1029").Value = 2.1
$ws.Range("O1029").Value = 3.75
$ws.Range("P1029").Value = 2.875
$ws.Range("Q1029").Value = -0.25
$ws.Range("R1029").Value = 1.925
$ws.Range("S1029").Value = 1.875
$ws.Range("T1029").Value = 2.75
$ws.Range("U1029").Value = 1.875
$ws.Range("V1029").Value = 1.925
$ws.Range("W1029").Value = 1.1
$ws.Range("X1029").Value = -1
$ws.Range("Y1029").Value = -1
$ws.Range("Z1029").Value = 0.925
$ws.Range("AA1029").Value = -1
$ws.Range("AB1029").Value = -1
$ws.Range("AC1029").Value = 0.925

# Row 1032 <- original row 1033
$ws.Range("B1032").Value = 6881621
$ws.Range("F1032").Value = 'Dynamo Dresden'
$ws.Range("G1032").Value = 'Hallescher FC'
$ws.Range("H1032").Value = 2
$ws.Range("I1032").Value = 1
$ws.Range("J1032").Value = 'H'
$ws.Range("K1032").Value = 1.444
$ws.Range("L1032").Value = 4.2
$ws.Range("M1032").Value = 6
$ws.Range("N1032").Value = 1.444
$ws.Range("O1032").Value = 4.333
$ws.Range("P1032").Value = 6
$ws.Range("Q1032").Value = -1.25
$ws.Range("R1032").Value = 1.95
$ws.Range("S1032").Value = 1.85
$ws.Range("T1032").Value = 3.25
$ws.Range("U1032").Value = 2
$ws.Range("V1032").Value = 1.8
$ws.Range("W1032").Value = 0.444
$ws.Range("X1032").Value = -1
$ws.Range("Y1032").Value = -1
$ws.Range("Z1032").Value = -0.5
$ws.Range("AA1032").Value = 0.425
$ws.Range("AB1032").Value = -0.5
$ws.Range("AC1032").Value = 0.4

# Row 1033 <- original row 1032
$ws.Range("B1033").Value = 6881622
$ws.Range("F1033").Value = 'Unterhaching'
$ws.Range("G1033").Value = 'RotWeiss Essen'
$ws.Range("H1033").Value = 4
$ws.Range("I1033").Value = 0
$ws.Range("J1033").Value = 'H'
$ws.Range("K1033").Value = 2.3
$ws.Range("L1033").Value = 3.5
$ws.Range("M1033").Value = 2.7
$ws.Range("N1033").Value = 2.375
$ws.Range("O1033").Value = 3.4
$ws.Range("P1033").Value = 2.7
$ws.Range("Q1033").Value = 0
$ws.Range("R1033").Value = 1.8
$ws.Range("S1033").Value = 2
$ws.Range("T1033").Value = 2.5
$ws.Range("U1033").Value = 1.875
$ws.Range("V1033").Value = 1.925
$ws.Range("W1033").Value = 1.375
$ws.Range("X1033").Value = -1
$ws.Range("Y1033").Value = -1
$ws.Range("Z1033").Value = 0.8
$ws.Range("AA1033").Value = -1
$ws.Range("AB1033").Value = 0.875
$ws.Range("AC1033").Value = -1

# Row 1149 <- original row 1153
$ws.Range("B1149").Value = 6880344
$ws.Range("F1149").Value = 'Waldhof Mannheim'
$ws.Range("G1149").Value = 'Dynamo Dresden'
$ws.Range("H1149").Value = 0
$ws.Range("I1149").Value = 2
$ws.Range("J1149").Value = 'A'
$ws.Range("K1149").Value = 4.8
$ws.Range("L1149").Value = 3.75
$ws.Range("M1149").Value = 1.615
$ws.Range("N1149").Value = 4
$ws.Range("O1149").Value = 3.4
$ws.Range("P1149").Value = 1.833
$ws.Range("Q1149").Value = 0.5
$ws.Range("R1149").Value = 1.9
$ws.Range("S1149").Value = 1.9
$ws.Range("T1149").Value = 2.5
$ws.Range("U1149").Value = 1.925
$ws.Range("V1149").Value = 1.875
$ws.Range("W1149").Value = -1
$ws.Range("X1149").Value = -1
$ws.Range("Y1149").Value = 0.833
$ws.Range("Z1149").Value = -1
$ws.Range("AA1149").Value = 0.8999999999999999
$ws.Range("AB1149").Value = -1
$ws.Range("AC1149").Value = 0.875

# Row 1150 <- original row 1151
$ws.Range("B1150").Value = 6881688
$ws.Range("F1150").Value = 'Vfb Lubeck'
$ws.Range("G1150").Value = '1860 Munich'
$ws.Range("H1150").Value = 1
$ws.Range("I1150").Value = 1
$ws.Range("J1150").Value = 'D'
$ws.Range("K1150").Value = 2.625
$ws.Range("L1150").Value = 3.4
$ws.Range("M1150").Value = 2.4
$ws.Range("N1150").Value = 2.9
$ws.Range("O1150").Value = 3.5
$ws.Range("P1150").Value = 2.15
$ws.Range("Q1150").Value = 0.25
$ws.Range("R1150").Value = 1.825
$ws.Range("S1150").Value = 1.975
$ws.Range("T1150").Value = 2.5
$ws.Range("U1150").Value = 1.8
$ws.Range("V1150").Value = 2
$ws.Range("W1150").Value = -1
$ws.Range("X1150").Value = 2.5
$ws.Range("Y1150").Value = -1
$ws.Range("Z1150").Value = 0.4125
$ws.Range("AA1150").Value = -0.5
$ws.Range("AB1150").Value = -1
$ws.Range("AC1150").Value = 1

# Row 1151 <- original row 1152
$ws.Range("B1151").Value = 6881685
$ws.Range("F1151").Value = 'Verl'
$ws.Range("G1151").Value = 'Saarbrucken'
$ws.Range("H1151").Value = 0
$ws.Range("I1151").Value = 0
$ws.Range("J1151").Value = 'D'
$ws.Range("K1151").Value = 2.75
$ws.Range("L1151").Value = 3.5
$ws.Range("M1151").Value = 2.25
$ws.Range("N1151").Value = 3
$ws.Range("O1151").Value = 3.4
$ws.Range("P1151").Value = 2.15
$ws.Range("Q1151").Value = 0.25
$ws.Range("R1151").Value = 1.875
$ws.Range("S1151").Value = 1.925
$ws.Range("T1151").Value = 2.5
$ws.Range("U1151").Value = 1.875
$ws.Range("V1151").Value = 1.925
$ws.Range("W1151").Value = -1
$ws.Range("X1151").Value = 2.4
$ws.Range("Y1151").Value = -1
$ws.Range("Z1151").Value = 0.4375
$ws.Range("AA1151").Value = -0.5
$ws.Range("AB1151").Value = -1
$ws.Range("AC1151").Value = 0.925

# Row 1152 <- original row 1149
$ws.Range("B1152").Value = 6881684
$ws.Range("F1152").Value = 'RotWeiss Essen'
$ws.Range("G1152").Value = 'FC Viktoria Kln'
$ws.Range("H1152").Value = 3
$ws.Range("I1152").Value = 1
$ws.Range("J1152").Value = 'H'
$ws.Range("K1152").Value = 2
$ws.Range("L1152").Value = 3.6
$ws.Range("M1152").Value = 3.2
$ws.Range("N1152").Value = 2.15
$ws.Range("O1152").Value = 3.5
$ws.Range("P1152").Value = 2.9
$ws.Range("Q1152").Value = -0.25
$ws.Range("R1152").Value = 1.95
$ws.Range("S1152").Value = 1.85
$ws.Range("T1152").Value = 2.75
$ws.Range("U1152").Value = 1.9
$ws.Range("V1152").Value = 1.9
$ws.Range("W1152").Value = 1.15
$ws.Range("X1152").Value = -1
$ws.Range("Y1152").Value = -1
$ws.Range("Z1152").Value = 0.95
$ws.Range("AA1152").Value = -1
$ws.Range("AB1152").Value = 0.8999999999999999
$ws.Range("AC1152").Value = -1

# Row 1153 <- original row 1150
$ws.Range("B1153").Value = 6880524
$ws.Range("F1153").Value = 'MSV Duisburg'
$ws.Range("G1153").Value = 'Hallescher FC'
$ws.Range("H1153").Value = 2
$ws.Range("I1153").Value = 3
$ws.Range("J1153").Value = 'A'
$ws.Range("K1153").Value = 2.2
$ws.Range("L1153").Value = 3.4
$ws.Range("M1153").Value = 2.9
$ws.Range("N1153").Value = 2.05
$ws.Range("O1153").Value = 3.6
$ws.Range("P1153").Value = 3.1
$ws.Range("Q1153").Value = -0.25
$ws.Range("R1153").Value = 1.825
$ws.Range("S1153").Value = 2.025
$ws.Range("T1153").Value = 2.75
$ws.Range("U1153").Value = 1.825
$ws.Range("V1153").Value = 2.025
$ws.Range("W1153").Value = -1
$ws.Range("X1153").Value = -1
$ws.Range("Y1153").Value = 2.1
$ws.Range("Z1153").Value = -1
$ws.Range("AA1153").Value = 1.025
$ws.Range("AB1153").Value = 0.825
$ws.Range("AC1153").Value = -1

# Row 1154 <- original row 1158
$ws.Range("B1154").Value = 6880480
$ws.Range("F1154").Value = 'Sandhausen'
$ws.Range("G1154").Value = 'Erzgebirge Aue'
$ws.Range("H1154").Value = 1
$ws.Range("I1154").Value = 0
$ws.Range("J1154").Value = 'H'
$ws.Range("K1154").Value = 2
$ws.Range("L1154").Value = 3.5
$ws.Range("M1154").Value = 3.25
$ws.Range("N1154").Value = 2
$ws.Range("O1154").Value = 3.5
$ws.Range("P1154").Value = 3.3
$ws.Range("Q1154").Value = -0.5
$ws.Range("R1154").Value = 2.05
$ws.Range("S1154").Value = 1.8
$ws.Range("T1154").Value = 2.75
$ws.Range("U1154").Value = 1.95
$ws.Range("V1154").Value = 1.9
$ws.Range("W1154").Value = 1
$ws.Range("X1154").Value = -1
$ws.Range("Y1154").Value = -1
$ws.Range("Z1154").Value = 1.05
$ws.Range("AA1154").Value = -1
$ws.Range("AB1154").Value = -1
$ws.Range("AC1154").Value = 0.8999999999999999

# Row 1155 <- original row 1157
$ws.Range("B1155").Value = 6881687
$ws.Range("F1155").Value = 'Unterhaching'
$ws.Range("G1155").Value = 'Freiburg II'
$ws.Range("H1155").Value = 1
$ws.Range("I1155").Value = 0
$ws.Range("J1155").Value = 'H'
$ws.Range("K1155").Value = 1.666
$ws.Range("L1155").Value = 3.8
$ws.Range("M1155").Value = 4.333
$ws.Range("N1155").Value = 1.7
$ws.Range("O1155").Value = 3.75
$ws.Range("P1155").Value = 4.2
$ws.Range("Q1155").Value = -0.75
$ws.Range("R1155").Value = 1.925
$ws.Range("S1155").Value = 1.875
$ws.Range("T1155").Value = 2.5
$ws.Range("U1155").Value = 1.95
$ws.Range("V1155").Value = 1.85
$ws.Range("W1155").Value = 0.7
$ws.Range("X1155").Value = -1
$ws.Range("Y1155").Value = -1
$ws.Range("Z1155").Value = 0.4625
$ws.Range("AA1155").Value = -0.5
$ws.Range("AB1155").Value = -1
$ws.Range("AC1155").Value = 0.8500000000000001

# Row 1157 <- original row 1155
$ws.Range("B1157").Value = 6881342
$ws.Range("F1157").Value = 'Arminia Bielefeld'
$ws.Range("G1157").Value = 'SSV Ulm 1846'
$ws.Range("H1157").Value = 0
$ws.Range("I1157").Value = 2
$ws.Range("J1157").Value = 'A'
$ws.Range("K1157").Value = 2.3
$ws.Range("L1157").Value = 3.5
$ws.Range("M1157").Value = 2.7
$ws.Range("N1157").Value = 2.2
$ws.Range("O1157").Value = 3.5
$ws.Range("P1157").Value = 2.8
$ws.Range("Q1157").Value = -0.25
$ws.Range("R1157").Value = 2
$ws.Range("S1157").Value = 1.85
$ws.Range("T1157").Value = 2.5
$ws.Range("U1157").Value = 1.825
$ws.Range("V1157").Value = 2.025
$ws.Range("W1157").Value = -1
$ws.Range("X1157").Value = -1
$ws.Range("Y1157").Value = 1.8
$ws.Range("Z1157").Value = -1
$ws.Range("AA1157").Value = 0.8500000000000001
$ws.Range("AB1157").Value = -1
$ws.Range("AC1157").Value = 1.025

# Row 1158 <- original row 1154
$ws.Range("B1158").Value = 6880370
$ws.Range("F1158").Value = 'FC Ingolstadt'
$ws.Range("G1158").Value = 'SC Preussen Munster'
$ws.Range("H1158").Value = 1
$ws.Range("I1158").Value = 1
$ws.Range("J1158").Value = 'D'
$ws.Range("K1158").Value = 2
$ws.Range("L1158").Value = 3.5
$ws.Range("M1158").Value = 3.25
$ws.Range("N1158").Value = 1.75
$ws.Range("O1158").Value = 4
$ws.Range("P1158").Value = 3.75
$ws.Range("Q1158").Value = -0.75
$ws.Range("R1158").Value = 1.975
$ws.Range("S1158").Value = 1.825
$ws.Range("T1158").Value = 3
$ws.Range("U1158").Value = 1.85
$ws.Range("V1158").Value = 1.95
$ws.Range("W1158").Value = -1
$ws.Range("X1158").Value = 3
$ws.Range("Y1158").Value = -1
$ws.Range("Z1158").Value = -1
$ws.Range("AA1158").Value = 0.825
$ws.Range("AB1158").Value = -1
$ws.Range("AC1158").Value = 0.95

# Row 1161 <- original row 1162
$ws.Range("B1161").Value = 6881691
$ws.Range("F1161").Value = 'Jahn Regensburg'
$ws.Range("G1161").Value = 'Arminia Bielefeld'
$ws.Range("H1161").Value = 2
$ws.Range("I1161").Value = 0
$ws.Range("J1161").Value = 'H'
$ws.Range("K1161").Value = 1.909
$ws.Range("L1161").Value = 3.6
$ws.Range("M1161").Value = 3.5
$ws.Range("N1161").Value = 2.25
$ws.Range("O1161").Value = 3.4
$ws.Range("P1161").Value = 2.875
$ws.Range("Q1161").Value = -0.25
$ws.Range("R1161").Value = 2.05
$ws.Range("S1161").Value = 1.8
$ws.Range("T1161").Value = 2.75
$ws.Range("U1161").Value = 2
$ws.Range("V1161").Value = 1.85
$ws.Range("W1161").Value = 1.25
$ws.Range("X1161").Value = -1
$ws.Range("Y1161").Value = -1
$ws.Range("Z1161").Value = 1.05
$ws.Range("AA1161").Value = -1
$ws.Range("AB1161").Value = -1
$ws.Range("AC1161").Value = 0.8500000000000001

# Row 1162 <- original row 1161
$ws.Range("B1162").Value = 6881692
$ws.Range("F1162").Value = 'Freiburg II'
$ws.Range("G1162").Value = 'Verl'
$ws.Range("H1162").Value = 0
$ws.Range("I1162").Value = 1
$ws.Range("J1162").Value = 'A'
$ws.Range("K1162").Value = 3.1
$ws.Range("L1162").Value = 3.5
$ws.Range("M1162").Value = 2
$ws.Range("N1162").Value = 3.3
$ws.Range("O1162").Value = 3.4
$ws.Range("P1162").Value = 2
$ws.Range("Q1162").Value = 0.25
$ws.Range("R1162").Value = 2
$ws.Range("S1162").Value = 1.8
$ws.Range("T1162").Value = 2.5
$ws.Range("U1162").Value = 1.95
$ws.Range("V1162").Value = 1.85
$ws.Range("W1162").Value = -1
$ws.Range("X1162").Value = -1
$ws.Range("Y1162").Value = 1
$ws.Range("Z1162").Value = -1
$ws.Range("AA1162").Value = 0.8
$ws.Range("AB1162").Value = -1
$ws.Range("AC1162").Value = 0.8500000000000001

# Row 1173 <- original row 1174
$ws.Range("B1173").Value = 6881696
$ws.Range("F1173").Value = 'Unterhaching'
$ws.Range("G1173").Value = 'Verl'
$ws.Range("K1173").Value = 2.3
$ws.Range("L1173").Value = 3.3
$ws.Range("M1173").Value = 2.75
$ws.Range("N1173").Value = 2.5
$ws.Range("O1173").Value = 3.3
$ws.Range("P1173").Value = 2.6
$ws.Range("Q1173").Value = 0
$ws.Range("R1173").Value = 1.875
$ws.Range("S1173").Value = 1.975
$ws.Range("T1173").Value = 2.75
$ws.Range("U1173").Value = 2.025
$ws.Range("V1173").Value = 1.825
$ws.Range("W1173").Value = 0
$ws.Range("X1173").Value = 0
$ws.Range("Y1173").Value = 0
$ws.Range("Z1173").Value = 0
$ws.Range("AA1173").Value = 0

# Row 1174 <- original row 1173
$ws.Range("B1174").Value = 6881694
$ws.Range("F1174").Value = 'RotWeiss Essen'
$ws.Range("G1174").Value = 'Freiburg II'
$ws.Range("K1174").Value = 1.55
$ws.Range("L1174").Value = 3.9
$ws.Range("M1174").Value = 5.5
$ws.Range("N1174").Value = 1.55
$ws.Range("O1174").Value = 3.9
$ws.Range("P1174").Value = 5.5
$ws.Range("Q1174").Value = -1
$ws.Range("R1174").Value = 1.95
$ws.Range("S1174").Value = 1.9
$ws.Range("T1174").Value = 2.5
$ws.Range("U1174").Value = 1.875
$ws.Range("V1174").Value = 1.975
$ws.Range("W1174").Value = 0
$ws.Range("X1174").Value = 0
$ws.Range("Y1174").Value = 0
$ws.Range("Z1174").Value = 0
$ws.Range("AA1174").Value = 0

# Standalone odds update for row 1169 (not part of a row shuffle)
$ws.Range("R1169").Value = 1.825
$ws.Range("S1169").Value = 2.025
